$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.6
$ws.Range("L2").Value = 5
$ws.Range("AG2").Value = 301

# --- Row 6 updates ---
$ws.Range("G6").Value = 2.4
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 3.2
$ws.Range("L6").Value = 4.33
$ws.Range("U6").Value = 2.38
$ws.Range("V6").Value = 1.53
$ws.Range("W6").Value = 5.5
$ws.Range("AI6").Value = 15
$ws.Range("AO6").Value = 15
$ws.Range("AU6").Value = 10
$ws.Range("AZ6").Value = 81
$ws.Range("BB6").Value = 450

# --- Row 8 is overwritten cell-by-cell with the data that used to be in row 9 ---
$ws.Cells.Item(8,1).Value = "ADYaA6BG"
$ws.Cells.Item(8,2).Value = "18/11/2024"
$ws.Cells.Item(8,3).Value = "21:30"
$ws.Cells.Item(8,4).Value = "URUGUAY - PRIMERA DIVISION"
$ws.Cells.Item(8,5).Value = "Cerro Largo"
$ws.Cells.Item(8,6).Value = "Wanderers"
$ws.Cells.Item(8,7).Value = 2.15
$ws.Cells.Item(8,8).Value = 3
$ws.Cells.Item(8,9).Value = 3.7
$ws.Cells.Item(8,10).Value = 2.88
$ws.Cells.Item(8,11).Value = 2.05
$ws.Cells.Item(8,12).Value = 4
$ws.Cells.Item(8,13).Value = 1.08
$ws.Cells.Item(8,14).Value = 8
$ws.Cells.Item(8,15).Value = 1.36
$ws.Cells.Item(8,16).Value = 3
$ws.Cells.Item(8,17).Value = 2.25
$ws.Cells.Item(8,18).Value = 1.62
$ws.Cells.Item(8,19).Value = 1.5
$ws.Cells.Item(8,20).Value = 2.5
$ws.Cells.Item(8,21).Value = 1.83
$ws.Cells.Item(8,22).Value = 1.83
$ws.Cells.Item(8,23).Value = 7
$ws.Cells.Item(8,24).Value = 9.5
$ws.Cells.Item(8,25).Value = 9.5
$ws.Cells.Item(8,26).Value = 19
$ws.Cells.Item(8,27).Value = 19
$ws.Cells.Item(8,28).Value = 34
$ws.Cells.Item(8,29).Value = 7.5
$ws.Cells.Item(8,30).Value = 6
$ws.Cells.Item(8,31).Value = 15
$ws.Cells.Item(8,32).Value = 51
$ws.Cells.Item(8,33).Value = 301
$ws.Cells.Item(8,34).Value = 10
$ws.Cells.Item(8,35).Value = 17
$ws.Cells.Item(8,36).Value = 13
$ws.Cells.Item(8,37).Value = 41
$ws.Cells.Item(8,38).Value = 34
$ws.Cells.Item(8,39).Value = 41
$ws.Cells.Item(8,40).Value = 4
$ws.Cells.Item(8,41).Value = 12
$ws.Cells.Item(8,42).Value = 23
$ws.Cells.Item(8,43).Value = 41
$ws.Cells.Item(8,44).Value = 67
$ws.Cells.Item(8,45).Value = 201
$ws.Cells.Item(8,46).Value = 2.5
$ws.Cells.Item(8,47).Value = 8.5
$ws.Cells.Item(8,48).Value = 67
$ws.Cells.Item(8,49).Value = 5.5
$ws.Cells.Item(8,50).Value = 21
$ws.Cells.Item(8,51).Value = 29
$ws.Cells.Item(8,52).Value = 67
$ws.Cells.Item(8,53).Value = 101
$ws.Cells.Item(8,54).Value = 251
$ws.Cells.Item(8,55).Value = 51
$ws.Cells.Item(8,56).Value = 51

# --- Row 9 is removed entirely (the sheet now ends at row 8) ---
$ws.Rows(9).Delete()

Write-Output "edit complete"
